$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -20.65879999999997
$ws.Range("A6").Value = -22.71320000000001
$ws.Range("A7").Value = -22.01550000000002
$ws.Range("B7").Value = 4.868900000000002
$ws.Range("B12").Value = 4.641
$ws.Range("E12").Value = 18.39320000000003
$ws.Range("D13").Value = -8.511299999999999
$ws.Range("D14").Value = -8.458000000000002
$ws.Range("B15").Value = 5.236499999999999
$ws.Range("A16").Value = -21.53519999999999
$ws.Range("D16").Value = -9.048100000000003
$ws.Range("D19").Value = -7.573699999999999
$ws.Range("A20").Value = -23.14780000000001
$ws.Range("B20").Value = 5.095399999999998
$ws.Range("B21").Value = 10.0448
$ws.Range("B22").Value = 8.818900000000003
$ws.Range("D22").Value = -8.181899999999999
$ws.Range("E22").Value = 16.73449999999999
$ws.Range("B23").Value = 9.0167
$ws.Range("A28").Value = -22.20979999999999
$ws.Range("A29").Value = -21.87389999999999
$ws.Range("B29").Value = 5.733699999999996
$ws.Range("E29").Value = 17.32090000000001
$ws.Range("A32").Value = -21.09409999999999
$ws.Range("B34").Value = 9.55150000000001
$ws.Range("E34").Value = 17.3256
$ws.Range("D36").Value = -8.223599999999999
$ws.Range("A40").Value = -19.1615
$ws.Range("B42").Value = 9.760499999999995
$ws.Range("B43").Value = 5.957299999999996
$ws.Range("E43").Value = 17.41680000000001
$ws.Range("B44").Value = 5.3717
$ws.Range("B45").Value = 5.124200000000001
$ws.Range("A46").Value = -22.0399
$ws.Range("B46").Value = 5.670000000000003
$ws.Range("D46").Value = -7.9853
$ws.Range("E48").Value = 17.5662
$ws.Range("B50").Value = 4.680599999999994
$ws.Range("D50").Value = -8.245299999999999
$ws.Range("A51").Value = -22.14589999999999
$ws.Range("B51").Value = 5.736599999999997
$ws.Range("A52").Value = -22.20499999999999
$ws.Range("A57").Value = -22.70800000000001
$ws.Range("A59").Value = -22.03550000000001
$ws.Range("E60").Value = 15.7085
$ws.Range("A62").Value = -22.06390000000001
$ws.Range("A66").Value = -21.5694
$ws.Range("B66").Value = 5.059399999999997
$ws.Range("B67").Value = 5.249299999999998
$ws.Range("E68").Value = 17.71450000000002
$ws.Range("E70").Value = 18.59860000000002
$ws.Range("A73").Value = -20.34469999999999
$ws.Range("E73").Value = 17.25140000000001
$ws.Range("A74").Value = -21.98889999999998
$ws.Range("B79").Value = 9.904200000000008
$ws.Range("B84").Value = 5.583799999999999
$ws.Range("E87").Value = 16.3055
$ws.Range("A92").Value = -21.46890000000002
$ws.Range("B92").Value = 4.608399999999993
$ws.Range("E92").Value = 19.10490000000002
$ws.Range("D95").Value = -8.268300000000004
$ws.Range("B97").Value = 5.758100000000001
$ws.Range("D97").Value = -8.451000000000001
$ws.Range("A100").Value = -21.98569999999999
$ws.Range("E101").Value = 16.89210000000001
